$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D7").Value = "No"
$ws.Range("D8").Value = "NULL"
